# Remove the "External Study ID" column from the Study sheet.
# This deletes column A (external_study_id) and shifts the remaining
# columns (parent_study_id, funding_source, principal_investigator,
# study_title, id, external_id) one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Study")

$ws.Columns.Item(1).Delete()
